# Update faturamento_anual data (ADD - Dados ADD PF)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - 2020
$ws.Range("B2").Value = 115509.02
$ws.Range("F2").Value = 447.7093798449612

# Row 3 - 2021
$ws.Range("B3").Value = 956160.54
$ws.Range("C3").Value = 727.7799777021743
$ws.Range("F3").Value = 442.8719499768412
$ws.Range("G3").Value = -1.080484369077828

# Row 4 - 2022
$ws.Range("B4").Value = 1772200.4
$ws.Range("C4").Value = 85.34548602057974
$ws.Range("F4").Value = 654.6732175840414
$ws.Range("G4").Value = 47.82449365291159

# Row 5 - 2023
$ws.Range("B5").Value = 2843142.44
$ws.Range("C5").Value = 60.43007551516184
$ws.Range("F5").Value = 852.772177564487
$ws.Range("G5").Value = 30.25921248336625

# Row 6 - 2024
$ws.Range("B6").Value = 4443078.12
$ws.Range("C6").Value = 56.27349715197527
$ws.Range("F6").Value = 938.1499408783784
$ws.Range("G6").Value = 10.01179043595557

# Row 7 - 2025
$ws.Range("B7").Value = 2281119.59
$ws.Range("C7").Value = -48.65902582869735
$ws.Range("D7").Value = 2306
$ws.Range("E7").Value = 2306
$ws.Range("F7").Value = 989.2105767562879
$ws.Range("G7").Value = 5.442694568642414
